# Insert a new data row at row 492 (pushing existing rows 492:559 down to 493:560)
# and populate it with the new "Femacal de La Calera - Zanahoria" record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(492).Insert()

$ws.Cells.Item(492, 1).Value  = 3
$ws.Cells.Item(492, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(492, 3).Value  = "Coquimbo"
$ws.Cells.Item(492, 4).Value  = 45077
$ws.Cells.Item(492, 5).Value  = 5
$ws.Cells.Item(492, 6).Value  = 100114013
$ws.Cells.Item(492, 7).Value  = "Zanahoria"
$ws.Cells.Item(492, 8).Value  = "Sin especificar"
$ws.Cells.Item(492, 9).Value  = "Primera"
$ws.Cells.Item(492, 10).Value = 250
$ws.Cells.Item(492, 11).Value = 7000
$ws.Cells.Item(492, 12).Value = 7500
$ws.Cells.Item(492, 13).Value = 7260
$ws.Cells.Item(492, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(492, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(492, 16).Value = 363
$ws.Cells.Item(492, 17).Value = 20
$ws.Cells.Item(492, 18).Value = "Hortaliza"
